# Generate Report for Handback
# Updates the localization-status report after a successful handback:
#  - Overview/zh-cn/de-de "Status" cells move from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - Latest Handback DateTime for zh-cn/de-de refreshed to the new handback time
#  - Error Detail cleared now that the handback file is up to date
#  - A few status/error columns are widened (or narrowed) to fit the new text

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-31 18:57:00"
$wsZhCn.Range("P2").ClearContents()
$wsZhCn.Columns.Item(3).ColumnWidth = 29.15
$wsZhCn.Columns.Item(16).ColumnWidth = 12.83

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-31 18:57:17"
$wsDeDe.Range("P2").ClearContents()
$wsDeDe.Columns.Item(3).ColumnWidth = 29.15
$wsDeDe.Columns.Item(16).ColumnWidth = 12.83
